# "fix in diagram IP->TP"
#
# 1) On the slide diagram, the rounded-rectangle box paired with
#    "Device TP-1" was mislabeled "KNX IP" -- it should read "KNX TP"
#    (it's a KNX TP device, matching its neighbour labelled "Device TP-1").
# 2) The cached text of every "datetimeFigureOut" date field (one per
#    slide layout, plus the slide master) advanced by a day, 07/11/2022
#    -> 08/11/2022, as PowerPoint re-stamped them on save.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Diagram label fix: "KNX IP" -> "KNX TP" on slide 1.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        $tr = $sh.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($j = 1; $j -le $paraCount; $j++) {
            $para = $tr.Paragraphs($j)
            $ptext = $para.Text.Trim()
            $nextText = ""
            if ($j + 1 -le $paraCount) {
                $nextText = $tr.Paragraphs($j + 1).Text.Trim()
            }
            if ($ptext -eq "KNX IP" -and $nextText -eq "Device TP-1") {
                $full = $tr.Characters($para.Start, $para.Length)
                $full.Text = "KNX TP"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Date placeholder cache refresh: 07/11/2022 -> 08/11/2022
#    (slide master + every custom layout).
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shape = $shapes.Item($k)
        if ($shape.Name -like "Date Placeholder*" -and $shape.HasTextFrame) {
            $dtr = $shape.TextFrame.TextRange
            if ($dtr.Text -eq "07/11/2022") {
                $dtr.Text = "08/11/2022"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
